$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column C (3rd column) width from 23.7109375 to 21.7109375
$ws.Columns.Item(3).ColumnWidth = 21.7109375

# Update data rows 2-22 for columns A, B, C, H, I, J
$data = @(
    @{Row=2; A=300; B=76657; C=9210; H=-22039; I=8912; J=456}
    @{Row=3; A=46; B=2868; C=646; H=-74862; I=-104715; J=277}
    @{Row=4; A=126; B=6710; C=898; H=92098; I=133769; J=45}
    @{Row=5; A=720; B=119922; C=9564; H=100146; I=187141; J=595}
    @{Row=6; A=959; B=276787; C=34271; H=82383; I=201657; J=620}
    @{Row=7; A=912; B=198551; C=58652; H=205232; I=404849; J=732}
    @{Row=8; A=453; B=59529; C=10614; H=8872; I=53474; J=513}
    @{Row=9; A=336; B=47835; C=26795; H=25643; I=56890; J=193}
    @{Row=10; A=237; B=23477; C=5100; H=-6824; I=21078; J=305}
    @{Row=11; A=529; B=69644; C=8732; H=14230; I=61096; J=534}
    @{Row=12; A=415; B=107979; C=61888; H=75887; I=105178; J=266}
    @{Row=13; A=825; B=188100; C=93465; H=-15419; I=17108; J=378}
    @{Row=14; A=1235; B=1317798; C=786198; H=668907; I=1130684; J=952}
    @{Row=15; A=517; B=200479; C=134647; H=97885; I=159838; J=747}
    @{Row=16; A=94; B=6074; C=1775; H=20496; I=46410; J=579}
    @{Row=17; A=122; B=10418; C=2338; H=-42198; I=-25957; J=402}
    @{Row=18; A=244; B=28917; C=18388; H=4821; I=-3909; J=500}
    @{Row=19; A=383; B=99483; C=81724; H=-12387; I=34227; J=452}
    @{Row=20; A=558; B=24722; C=12674; H=31670; I=51872; J=335}
    @{Row=21; A=543; B=31199; C=9198; H=50278; I=142983; J=592}
    @{Row=22; A=633; B=227985; C=185155; H=247112; I=442549; J=713}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = $item.I
    $ws.Cells.Item($r, 10).Value = $item.J
}
